$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, reusing the header style from G1 (avoids creating a
# duplicate style entry in styles.xml)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill H2:H8 with 0 values for the new "Save" column (plain numeric cells,
# same as the rest of the data rows - no special style)
$ws.Range("H2:H8").Value = 0
